$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Dyn_Rule")
$ws2 = $wb.Worksheets.Item("Field_Rule")

# Header H1: change from Condition_Value placeholder text... Actually just set value
$ws1.Range("H1").Value = "Condition_Value"

# Row2 data
$ws1.Range("F2").Value = "013"
$ws1.Range("G2").Value = "//option[@value=`"=`"] "
$ws1.Range("H2").Value = "1000"
$ws1.Range("I2").Value = "014"

$ws1.Range("A1:J2").Select()
